# Updates from the 6/18 status meeting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window / view position (best effort; engine may not persist these) ---
$win = $excel.ActiveWindow
$win.Left = 1860
$win.Top = 0

# --- Row 100: new action item row ---
$ws.Range("B100").Value = "Find out when UPT 5.0 will be officially released and required on the tech stack."
$ws.Range("C100").Value = "JJ Pan"
$ws.Range("D100").Value = 39981
$ws.Range("E100").Value = "Assigned"

# --- Row 101: new action item row ---
$ws.Range("B101").Value = "Create/update the test plan for testing the restore of caArray data from tape backup."
$ws.Range("C101").Value = "Mike Hunter, Marina Omelchenko, Winston Cheng"
$ws.Range("D101").Value = 39981
$ws.Range("E101").Value = "Assigned"

# Rows 100 and 101 grow taller (wrapped text) to match the new content.
$ws.Range("A100:E101").EntireRow.RowHeight = 31

# --- Extend the table down to row 122 with blank, numbered rows ---
for ($i = 103; $i -le 122; $i++) {
    $ws.Range("A$i").Value = $i - 1
}
$ws.Range("A102:E102").Copy()
$ws.Range("A103:E122").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A103:E122").EntireRow.RowHeight = 16

# --- Update the view / selection to where the meeting notes were added ---
$ws.Range("E102").Select()
